# Versie voor aanpassing na feedback
#
# 1. Append a period ("." as its own run) to the end of the paragraph
#    "Klasse Scoreboard staat in voor de score van de huidige game".
# 2. Insert a brand-new paragraph right after it describing the abstract
#    Enum class (with spell-check proofErr markers around the words Word's
#    checker would flag: "num" (of "Enum"), "Seriazable", "Comparable").

$d = $word.ActiveDocument

# --- Locate the target paragraph robustly (search by its text rather than
#     a hard-coded index) ----------------------------------------------------
$targetIndex = -1
$i = 1
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Klasse Scoreboard staat in voor de score van de huidige game*") {
        $targetIndex = $i
        break
    }
    $i++
}

if ($targetIndex -eq -1) {
    Write-Output "ERROR: target paragraph 'Klasse Scoreboard ...' not found"
} else {
    $target = $d.Paragraphs.Item($targetIndex)
    $r = $target.Range

    # Recover the paragraph's own opening <w:p ...> tag (carrying its
    # w14:paraId / w14:textId / rsid identity) so that rebuilding it below
    # keeps those attributes intact instead of losing them. Falls back to a
    # plain <w:p> if anything about that lookup doesn't pan out.
    $openTag = "<w:p>"
    try {
        $fullXml = $r.WordOpenXML
        if ($fullXml -match '(<w:p[ >][^>]*>)<w:r><w:t>Klasse Scoreboard') {
            $openTag = $matches[1]
        }
    } catch { }

    # --- Step 1: rewrite the paragraph so the trailing "." lives in its own
    #     run, just like typing the sentence and then appending a period as
    #     a separate edit afterwards. Rebuilding the paragraph via InsertXML
    #     (a WordProcessingML fragment) keeps that run boundary on save. ----
    $bodyParagraph = $openTag + "<w:r><w:t>Klasse Scoreboard staat in voor de score van de huidige game</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>"

    $xmlPeriod = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
$bodyParagraph
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $r.InsertXML($xmlPeriod)

    # --- Step 2: add a new paragraph right after it with the remark about
    #     the abstract Enum class / Serializable & Comparable interfaces. ---
    $target2 = $d.Paragraphs.Item($targetIndex)
    $target2.Range.InsertParagraphAfter()

    $newPara = $d.Paragraphs.Item($targetIndex + 1)

    $xmlNew = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>Abstracte klasse E</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>num</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> met Interface </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Seriazable</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> &amp; </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Comparable</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> wordt automatisch door Java aangemaakt.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
    $newPara.Range.InsertXML($xmlNew)

    Write-Output "OK: paragraph $targetIndex updated, new paragraph inserted at $($targetIndex + 1)"
}
